# Update "Notified Production Wind" data:
# - Shift every timestamp in column A forward by 8 days (Astro -> Dabaca location change)
# - Replace the corresponding "Notified Production (MW)" values in column B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 45807.01041666666, 340),
    @(3, 45807.02083333334, 343),
    @(4, 45807.03125, 346),
    @(5, 45807.04166666666, 348),
    @(6, 45807.05208333334, 444),
    @(7, 45807.0625, 448),
    @(8, 45807.07291666666, 451),
    @(9, 45807.08333333334, 455),
    @(10, 45807.09375, 627),
    @(11, 45807.10416666666, 632),
    @(12, 45807.11458333334, 636),
    @(13, 45807.125, 640),
    @(14, 45807.13541666666, 811),
    @(15, 45807.14583333334, 816),
    @(16, 45807.15625, 821),
    @(17, 45807.16666666666, 826),
    @(18, 45807.17708333334, 1017),
    @(19, 45807.1875, 1021),
    @(20, 45807.19791666666, 1026),
    @(21, 45807.20833333334, 1030),
    @(22, 45807.21875, 1152),
    @(23, 45807.22916666666, 1153),
    @(24, 45807.23958333334, 1155),
    @(25, 45807.25, 1158),
    @(26, 45807.26041666666, 1166),
    @(27, 45807.27083333334, 1166),
    @(28, 45807.28125, 1165),
    @(29, 45807.29166666666, 1165),
    @(30, 45807.30208333334, 1123),
    @(31, 45807.3125, 1124),
    @(32, 45807.32291666666, 1126),
    @(33, 45807.33333333334, 1128),
    @(34, 45807.34375, 1140),
    @(35, 45807.35416666666, 1140),
    @(36, 45807.36458333334, 1141),
    @(37, 45807.375, 1141),
    @(38, 45807.38541666666, 1188),
    @(39, 45807.39583333334, 1188),
    @(40, 45807.40625, 1188),
    @(41, 45807.41666666666, 1187),
    @(42, 45807.42708333334, 1158),
    @(43, 45807.4375, 1159),
    @(44, 45807.44791666666, 1160),
    @(45, 45807.45833333334, 1161),
    @(46, 45807.46875, 1152),
    @(47, 45807.47916666666, 1153),
    @(48, 45807.48958333334, 1153),
    @(49, 45807.5, 1153),
    @(50, 45807.51041666666, 1167),
    @(51, 45807.52083333334, 1166),
    @(52, 45807.53125, 1165),
    @(53, 45807.54166666666, 1164),
    @(54, 45807.55208333334, 1192),
    @(55, 45807.5625, 1191),
    @(56, 45807.57291666666, 1191),
    @(57, 45807.58333333334, 1190),
    @(58, 45807.59375, 1210),
    @(59, 45807.60416666666, 1210),
    @(60, 45807.61458333334, 1210),
    @(61, 45807.625, 1210),
    @(62, 45807.63541666666, 1214),
    @(63, 45807.64583333334, 1215),
    @(64, 45807.65625, 1216),
    @(65, 45807.66666666666, 1216),
    @(66, 45807.67708333334, 1176),
    @(67, 45807.6875, 1175),
    @(68, 45807.69791666666, 1174),
    @(69, 45807.70833333334, 1172),
    @(70, 45807.71875, 1114),
    @(71, 45807.72916666666, 1110),
    @(72, 45807.73958333334, 1105),
    @(73, 45807.75, 1100),
    @(74, 45807.76041666666, 982),
    @(75, 45807.77083333334, 978),
    @(76, 45807.78125, 973),
    @(77, 45807.79166666666, 969),
    @(78, 45807.80208333334, 907),
    @(79, 45807.8125, 906),
    @(80, 45807.82291666666, 905),
    @(81, 45807.83333333334, 905),
    @(82, 45807.84375, 912),
    @(83, 45807.85416666666, 913),
    @(84, 45807.86458333334, 915),
    @(85, 45807.875, 916),
    @(86, 45807.88541666666, 1008),
    @(87, 45807.89583333334, 1010),
    @(88, 45807.90625, 1012),
    @(89, 45807.91666666666, 1015),
    @(90, 45807.92708333334, 1054),
    @(91, 45807.9375, 1056),
    @(92, 45807.94791666666, 1058),
    @(93, 45807.95833333334, 1060),
    @(94, 45807.96875, 0),
    @(95, 45807.97916666666, 0),
    @(96, 45807.98958333334, 0),
    @(97, 45808, 0)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $tsValue = $entry[1]
    $mwValue = $entry[2]
    $ws.Cells.Item($row, 1).Value = $tsValue
    $ws.Cells.Item($row, 2).Value = $mwValue
}
